$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (row => column letter => value) reflecting updated NATMI computation
$data = @{
    2  = @{ E=3; G=2.748846; H=8.246537999999999; I=0.1900189725460084; J=0.1900189725460084; K=3; M=2.633247; N=7.899741; O=0.2470706468821852; P=0.2470706468821852; Q=7.238390482961999; R=65.14551434665799; S=0.04694811046683049; T=0.04694811046683049 }
    3  = @{ E=3; G=2.748846; H=8.246537999999999; I=0.1900189725460084; J=0.1900189725460084; K=3; M=5.113233999999999; N=15.339702; O=0.4797613106708119; P=0.479761310670812; Q=14.055492827964; R=126.499435451676; S=0.09116375132099402; T=0.09116375132099401 }
    4  = @{ E=3; G=2.748846; H=8.246537999999999; I=0.1900189725460084; J=0.1900189725460084; K=3; M=2.911389666666667; N=8.734169; O=0.2731680424470029; P=0.2731680424470029; Q=8.002961839657999; R=72.02665655692199; S=0.0519071107581839; T=0.05190711075818388 }
    5  = @{ E=3; G=7.625013333333333; H=22.87504; I=0.5270928961642867; J=0.5270928961642866; K=3; M=2.633247; N=7.899741; O=0.2470706468821852; P=0.2470706468821852; Q=20.07854348496; R=180.70689136464; S=0.1302291828223148; T=0.1302291828223148 }
    6  = @{ E=3; G=7.625013333333333; H=22.87504; I=0.5270928961642867; J=0.5270928961642866; K=3; M=5.113233999999999; N=15.339702; O=0.4797613106708119; P=0.479761310670812; Q=38.98847742645332; R=350.8962968380799; S=0.2528787787090523; T=0.2528787787090523 }
    7  = @{ E=3; G=7.625013333333333; H=22.87504; I=0.5270928961642867; J=0.5270928961642866; K=3; M=2.911389666666667; N=8.734169; O=0.2731680424470029; P=0.2731680424470029; Q=22.19938502686222; R=199.79446524176; S=0.1439849346329195; T=0.1439849346329195 }
    8  = @{ E=3; G=4.092306666666667; H=12.27692; I=0.282888131289705; J=0.282888131289705; K=3; M=2.633247; N=7.899741; O=0.2470706468821852; P=0.2470706468821852; Q=10.77605425308; R=96.98448827772; S=0.06989335359303996; T=0.06989335359303996 }
    9  = @{ E=3; G=4.092306666666667; H=12.27692; I=0.282888131289705; J=0.282888131289705; K=3; M=5.113233999999999; N=15.339702; O=0.4797613106708119; P=0.479761310670812; Q=20.92492158642666; R=188.32429427784; S=0.1357187806407656; T=0.1357187806407656 }
    10 = @{ E=3; G=4.092306666666667; H=12.27692; I=0.282888131289705; J=0.282888131289705; K=3; M=2.911389666666667; N=8.734169; O=0.2731680424470029; P=0.2731680424470029; Q=11.91429934216444; R=107.22869407948; S=0.07727599705589946; T=0.07727599705589945 }
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$r").Value = $rowVals[$col]
    }
}
